$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 30.555555
$ws.Range("I6").Value = 30.555555
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 91.66666499999999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 20.33333500000001
$ws.Range("N6").ClearContents()
$ws.Range("H18").Value = 914.3333
$ws.Range("I18").Value = 914.3333
$ws.Range("K18").Value = 914.3333
$ws.Range("M18").Value = -630.3333
$ws.Range("H40").Value = 2444.818
$ws.Range("J40").Value = 1874.125
$ws.Range("L40").Value = 1874.125
$ws.Range("N40").Value = -2224.125
$ws.Range("H51").Value = 12499.5
$ws.Range("H81").Value = 56000
$ws.Range("J81").Value = 56000
$ws.Range("L81").Value = 56000
$ws.Range("N81").Value = -57996
$ws.Range("H84").Value = 56000
$ws.Range("J84").Value = 56000
$ws.Range("L84").Value = 168000
$ws.Range("N84").Value = -177984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3476.5715
$ws.Range("I32").Value = 3476.5715
$ws.Range("K32").Value = 3476.5715
$ws.Range("M32").Value = -3189.5715
$ws.Range("H61").Value = 3572
$ws.Range("I61").Value = 3572
$ws.Range("K61").Value = 3572
$ws.Range("M61").Value = -3360
$ws.Range("H63").Value = 2000
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2000
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 10000
$ws.Range("N66").Value = -16864
$ws.Range("H92").Value = 72996.7
$ws.Range("J92").Value = 79218.55499999999
$ws.Range("L92").Value = 79218.55499999999
$ws.Range("N92").Value = -84210.55499999999
$ws.Range("H136").Value = 3572
$ws.Range("I136").Value = 3572
$ws.Range("K136").Value = 10716
$ws.Range("M136").Value = -8166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2333.3333
$ws.Range("I20").Value = 1750
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 1750
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -1503
$ws.Range("N20").Value = -3994
$ws.Range("H22").Value = 801.36365
$ws.Range("I22").Value = 757.3333
$ws.Range("K22").Value = 757.3333
$ws.Range("M22").Value = -584.3333
$ws.Range("H92").Value = 102645
$ws.Range("J92").Value = 102645
$ws.Range("L92").Value = 102645
$ws.Range("N92").Value = -107637
$ws.Range("H94").Value = 123656.22
$ws.Range("I94").Value = 123656.22
$ws.Range("K94").Value = 123656.22
$ws.Range("M94").Value = -123205.22
$ws.Range("H138").Value = 39000
$ws.Range("J138").Value = 39000
$ws.Range("L138").Value = 39000
$ws.Range("N138").Value = -49280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9826
$ws.Range("H32").Value = 2396.5
$ws.Range("I32").Value = 793
$ws.Range("K32").Value = 793
$ws.Range("M32").Value = -477
$ws.Range("H33").Value = 18118.428
$ws.Range("I33").Value = 1366
$ws.Range("K33").Value = 1366
$ws.Range("M33").Value = -987
$ws.Range("H35").Value = 2206.8572
$ws.Range("I35").Value = 2424.6667
$ws.Range("J35").Value = 900
$ws.Range("K35").Value = 2424.6667
$ws.Range("L35").Value = 900
$ws.Range("M35").Value = -2130.6667
$ws.Range("N35").Value = -1488
$ws.Range("H36").Value = 39999.4
$ws.Range("I36").Value = 19999
$ws.Range("K36").Value = 19999
$ws.Range("M36").Value = -19611
$ws.Range("H40").Value = 39999.4
$ws.Range("I40").Value = 19999
$ws.Range("K40").Value = 19999
$ws.Range("M40").Value = -19839
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 2306.7273
$ws.Range("I134").Value = 2263.7778
$ws.Range("K134").Value = 6791.3334
$ws.Range("M134").Value = -4256.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1804
$ws.Range("J63").Value = 2666.6667
$ws.Range("L63").Value = 8000.000100000001
$ws.Range("N63").Value = -9498.000100000001
$ws.Range("H66").Value = 1804
$ws.Range("J66").Value = 2666.6667
$ws.Range("L66").Value = 24000.0003
$ws.Range("N66").Value = -31488.0003
$ws.Range("H117").Value = 1874.3
$ws.Range("I117").Value = 94
$ws.Range("J117").Value = 2637.2856
$ws.Range("K117").Value = 282
$ws.Range("L117").Value = 7911.8568
$ws.Range("M117").Value = 3160
$ws.Range("N117").Value = -14795.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H9").Value = 303.5
$ws.Range("I9").Value = 303.5
$ws.Range("K9").Value = 303.5
$ws.Range("M9").Value = -133.5
$ws.Range("H31").Value = 2239.8
$ws.Range("I31").Value = 2239.8
$ws.Range("K31").Value = 2239.8
$ws.Range("M31").Value = -1947.8
$ws.Range("H37").Value = 2239.8
$ws.Range("I37").Value = 2239.8
$ws.Range("K37").Value = 2239.8
$ws.Range("M37").Value = -1962.8
$ws.Range("H41").Value = 22249.5
$ws.Range("I41").Value = 19666
$ws.Range("K41").Value = 19666
$ws.Range("M41").Value = -19311
$ws.Range("H80").Value = 4363.8887
$ws.Range("I80").Value = 3248.75
$ws.Range("K80").Value = 3248.75
$ws.Range("M80").Value = -2250.75
$ws.Range("H83").Value = 4363.8887
$ws.Range("I83").Value = 3248.75
$ws.Range("K83").Value = 16243.75
$ws.Range("M83").Value = -11251.75
$ws.Range("H105").Value = 11622
$ws.Range("J105").Value = 11622
$ws.Range("L105").Value = 11622
$ws.Range("N105").Value = -18610
$ws.Range("H113").Value = 863.2
$ws.Range("I113").Value = 863.2
$ws.Range("K113").Value = 863.2
$ws.Range("M113").Value = 1306.8
$ws.Range("H132").Value = 1600.5714
$ws.Range("I132").Value = 1600.5714
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4801.7142
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2271.7142
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 231.25
$ws.Range("J22").Value = 250
$ws.Range("L22").Value = 250
$ws.Range("N22").Value = -840
$ws.Range("H27").Value = 231.25
$ws.Range("J27").Value = 250
$ws.Range("L27").Value = 250
$ws.Range("N27").Value = -464
$ws.Range("H35").Value = 11016.111
$ws.Range("I35").Value = 1430
$ws.Range("K35").Value = 1430
$ws.Range("M35").Value = -1094
$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1312
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 1883
$ws.Range("I61").Value = 1824.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1824.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1622.5
$ws.Range("N61").Value = -2404
$ws.Range("H82").Value = 4598.8335
$ws.Range("I82").Value = 3538.6
$ws.Range("K82").Value = 3538.6
$ws.Range("M82").Value = -3177.6
$ws.Range("H85").Value = 4598.8335
$ws.Range("I85").Value = 3538.6
$ws.Range("K85").Value = 3538.6
$ws.Range("M85").Value = -2290.6
$ws.Range("H113").Value = 1883
$ws.Range("I113").Value = 1824.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1824.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 345.5
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 1201.3334
$ws.Range("I132").Value = 1201.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3604.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1074.0002
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4749.5
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 4499
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 4499
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -5747
$ws.Range("H65").Value = 4749.5
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 4499
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 22495
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -28735
$ws.Range("H105").Value = 26665.666
$ws.Range("J105").Value = 26665.666
$ws.Range("L105").Value = 26665.666
$ws.Range("N105").Value = -33653.666
$ws.Range("H132").Value = 1287.5
$ws.Range("I132").Value = 1283.3334
$ws.Range("K132").Value = 3850.0002
$ws.Range("M132").Value = -1320.0002
